$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write the new column K values (year 2020 data) and matching fonts ---
# "Total" rows (4, 30 and the republic sub-total rows) use bold 9pt Times New Roman,
# matching the style already used by the corresponding cells in column I/J.
# The remaining (urban/rural) rows use the regular (non-bold) 9pt Times New Roman font.

$boldRows = @(4, 5, 8, 11, 14, 17, 20, 23, 26, 29, 30)
$regRows = @(6, 7, 9, 10, 12, 13, 15, 16, 18, 19, 21, 22, 24, 25, 27, 28)

$values = @{
    4 = 2020
    5 = 22
    6 = 29.1
    7 = 20.2
    8 = 26.8
    9 = 39.8
    10 = 22.7
    11 = 22
    12 = 33.8
    13 = 18.8
    14 = 28
    15 = 38.7
    16 = 22
    17 = 35.1
    18 = 45.1
    19 = 33.3
    20 = 19.4
    21 = 13
    22 = 19.9
    23 = 26.2
    24 = 52.8
    25 = 22.5
    26 = 20.1
    27 = 33.6
    28 = 16.9
    29 = 23.5
    30 = 30.7
}

foreach ($r in $boldRows) {
    $cell = $ws.Range("K$r")
    $cell.Value = $values[$r]
    $cell.Font.Name = "Times New Roman"
    $cell.Font.Size = 9
    $cell.Font.Bold = $true
    $cell.VerticalAlignment = -4107
}

foreach ($r in $regRows) {
    $cell = $ws.Range("K$r")
    $cell.Value = $values[$r]
    $cell.Font.Name = "Times New Roman"
    $cell.Font.Size = 9
    $cell.Font.Bold = $false
    $cell.VerticalAlignment = -4107
}

# --- Apply borders ---
# K4: medium border on top and bottom (matches the header row border)
$ws.Range("K4").Borders.Item(9).Weight = -4138
$ws.Range("K4").Borders.Item(8).Weight = -4138

# K30: medium border on the bottom only (matches the footer row border)
$ws.Range("K30").Borders.Item(9).Weight = -4138

# --- Selection, matching the saved cursor position in the workbook ---
$ws.Range("K18").Select() | Out-Null
